# Auto-generated Excel COM-interop script
# Applies numeric cell updates to match target diff for Sheets/Coeurl_Profits.xlsx
$wb = $excel.ActiveWorkbook

# --- ALC row 58 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 791.7273
$ws.Range("I58").Value = 223.5
$ws.Range("J58").Value = 2307
$ws.Range("K58").Value = 670.5
$ws.Range("L58").Value = 6921
$ws.Range("M58").Value = -520.5
$ws.Range("N58").Value = -7221

# --- ALC row 76 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5842.7144
$ws.Range("I76").Value = 5780.4
$ws.Range("K76").Value = 5780.4
$ws.Range("M76").Value = -5465.4

# --- ALC row 79 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5842.7144
$ws.Range("I79").Value = 5780.4
$ws.Range("K79").Value = 5780.4
$ws.Range("M79").Value = -4688.4

# --- ALC row 109 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 42983.5
$ws.Range("J109").Value = 42983.5
$ws.Range("L109").Value = 42983.5
$ws.Range("N109").Value = -45757.5

# --- ALC row 131 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 25888.143
$ws.Range("I131").Value = 4242.4
$ws.Range("J131").Value = 80002.5
$ws.Range("K131").Value = 12727.2
$ws.Range("L131").Value = 240007.5
$ws.Range("M131").Value = -7687.199999999999
$ws.Range("N131").Value = -250087.5

# --- ALC row 135 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1344.2858
$ws.Range("I135").Value = 892.8182
$ws.Range("J135").Value = 2999.6667
$ws.Range("K135").Value = 8035.3638
$ws.Range("L135").Value = 26997.0003
$ws.Range("M135").Value = -5500.3638
$ws.Range("N135").Value = -32067.0003

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 28541.186
$ws.Range("I138").Value = 170028.17
$ws.Range("J138").Value = 5597.3516
$ws.Range("K138").Value = 510084.51
$ws.Range("L138").Value = 16792.0548
$ws.Range("M138").Value = -504944.51
$ws.Range("N138").Value = -27072.0548

# --- ARM row 45 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 74858
$ws.Range("I45").Value = 74858
$ws.Range("K45").Value = 74858
$ws.Range("M45").Value = -74481

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 15410.308
$ws.Range("I74").Value = 1654.875
$ws.Range("K74").Value = 1654.875
$ws.Range("M74").Value = -780.875

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 15410.308
$ws.Range("I77").Value = 1654.875
$ws.Range("K77").Value = 8274.375
$ws.Range("M77").Value = -3906.375

# --- ARM row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2386.7817
$ws.Range("I132").Value = 2162.976
$ws.Range("K132").Value = 6488.928
$ws.Range("M132").Value = -3958.928

# --- BSM row 19 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 540
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -327

# --- BSM row 105 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1817.8422
$ws.Range("I105").Value = 1467.1818
$ws.Range("J105").Value = 2300
$ws.Range("K105").Value = 1467.1818
$ws.Range("L105").Value = 2300
$ws.Range("M105").Value = 279.8181999999999
$ws.Range("N105").Value = -5794

# --- CRP row 22 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 609.8
$ws.Range("I22").Value = 483.33334
$ws.Range("J22").Value = 799.5
$ws.Range("K22").Value = 483.33334
$ws.Range("L22").Value = 799.5
$ws.Range("M22").Value = -133.33334
$ws.Range("N22").Value = -1499.5

# --- CRP row 60 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 7688.1113
$ws.Range("J60").Value = 19999
$ws.Range("L60").Value = 19999
$ws.Range("N60").Value = -21021

# --- CRP row 105 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2060
$ws.Range("I105").Value = 1497.5
$ws.Range("K105").Value = 1497.5
$ws.Range("M105").Value = 249.5

# --- CUL row 55 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 762.5
$ws.Range("I55").Value = 156.25
$ws.Range("J55").Value = 1975
$ws.Range("K55").Value = 468.75
$ws.Range("L55").Value = 5925
$ws.Range("M55").Value = -291.75
$ws.Range("N55").Value = -6279

# --- CUL row 113 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 970.4286
$ws.Range("I113").Value = 197.5
$ws.Range("K113").Value = 592.5
$ws.Range("M113").Value = 1577.5

# --- CUL row 129 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 625.9286
$ws.Range("I129").Value = 520.1539
$ws.Range("J129").Value = 2001
$ws.Range("K129").Value = 1560.4617
$ws.Range("L129").Value = 6003
$ws.Range("M129").Value = 3439.5383
$ws.Range("N129").Value = -16003

# --- CUL row 131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 40531.81
$ws.Range("I131").Value = 59707.53
$ws.Range("K131").Value = 179122.59
$ws.Range("M131").Value = -174082.59

# --- CUL row 137 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4039.8462
$ws.Range("I137").Value = 2391.4
$ws.Range("J137").Value = 5070.125
$ws.Range("K137").Value = 7174.200000000001
$ws.Range("L137").Value = 15210.375
$ws.Range("M137").Value = -2074.200000000001
$ws.Range("N137").Value = -25410.375

# --- CUL row 140 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3272.9167
$ws.Range("I140").Value = 3343.2727
$ws.Range("K140").Value = 10029.8181
$ws.Range("M140").Value = -4849.8181

# --- GSM row 18 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 14333.333
$ws.Range("J18").Value = 12500
$ws.Range("L18").Value = 12500
$ws.Range("N18").Value = -13086

# --- GSM row 26 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

# --- GSM row 49 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("N49").Value = -20368

# --- GSM row 50 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()

# --- LTW row 7 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21500.8
$ws.Range("I7").Value = 24876
$ws.Range("K7").Value = 24876
$ws.Range("M7").Value = -24764

# --- LTW row 16 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 55557452
$ws.Range("I16").Value = 55557452
$ws.Range("K16").Value = 55557452
$ws.Range("M16").Value = -55557282

# --- LTW row 55 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 679.06665
$ws.Range("I55").Value = 475.875
$ws.Range("J55").Value = 911.2857
$ws.Range("K55").Value = 475.875
$ws.Range("L55").Value = 911.2857
$ws.Range("M55").Value = -302.875
$ws.Range("N55").Value = -1257.2857

# --- LTW row 61 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5053
$ws.Range("I61").Value = 5099.231
$ws.Range("J61").Value = 4752.5
$ws.Range("K61").Value = 5099.231
$ws.Range("L61").Value = 4752.5
$ws.Range("M61").Value = -4897.231
$ws.Range("N61").Value = -5156.5

# --- LTW row 113 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5053
$ws.Range("I113").Value = 5099.231
$ws.Range("J113").Value = 4752.5
$ws.Range("K113").Value = 5099.231
$ws.Range("L113").Value = 4752.5
$ws.Range("M113").Value = -2929.231
$ws.Range("N113").Value = -9092.5

# --- LTW row 126 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 21500.8
$ws.Range("I126").Value = 24876
$ws.Range("K126").Value = 74628
$ws.Range("M126").Value = -72158

# --- LTW row 136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 53610.3
$ws.Range("I136").Value = 75057.71000000001
$ws.Range("K136").Value = 225173.13
$ws.Range("M136").Value = -222623.13

# --- WVR row 100 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1248.3478
$ws.Range("I100").Value = 466.0909
$ws.Range("J100").Value = 1965.4166
$ws.Range("K100").Value = 932.1818
$ws.Range("L100").Value = 3930.8332
$ws.Range("M100").Value = -391.1818
$ws.Range("N100").Value = -5012.8332

# --- WVR row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1329.638
$ws.Range("I132").Value = 1195.88
$ws.Range("J132").Value = 2165.625
$ws.Range("K132").Value = 3587.64
$ws.Range("L132").Value = 6496.875
$ws.Range("M132").Value = -1057.64
$ws.Range("N132").Value = -11556.875
